$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$r = $p1.Range
Get-Member -InputObject $r
